$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab) to reflect the new "through" date
$ws.Name = "Through 2022-02-12"

# Update the row label for February to the new "through" date
$ws.Range("A3").Value = "February (through 02-12)"

# Update February row values (row 3) for years 2015-2022 (columns B-I)
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 31
$ws.Range("E3").Value = 23
$ws.Range("F3").Value = 12
$ws.Range("G3").Value = 28
$ws.Range("H3").Value = 61
$ws.Range("I3").Value = 56

# Update Total row values (row 4) for years 2015-2022 (columns B-I)
$ws.Range("B4").Value = 31
$ws.Range("C4").Value = 63
$ws.Range("D4").Value = 106
$ws.Range("E4").Value = 109
$ws.Range("F4").Value = 61
$ws.Range("G4").Value = 102
$ws.Range("H4").Value = 278
$ws.Range("I4").Value = 217
